$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Ryan's "First Resume Interviews" count from 0 to 1 (cell F7)
$ws.Range("F7").Value = 1

# Refresh the chart so its cached series data picks up the new value
if ($ws.ChartObjects().Count -gt 0) {
    $ws.ChartObjects().Item(1).Chart.Refresh()
}

# Update the selected cell to match the final saved state
$ws.Range("G15").Select()
